$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS"): the cash-flow table's
#    style changes from {DF3FFDFF-BDF7-400A-9402-3547228F21C4} to
#    {FF854312-E763-4362-AE32-1D3A4BC16B2E}.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FF854312-E763-4362-AE32-1D3A4BC16B2E}")
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme switches from the custom "Integral" palette back to the
#    stock Office palette. Push every Office theme color into the live
#    ThemeColorScheme (order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#    COM RGB values are encoded 0x00BBGGRR, so convert each RRGGBB hex first.
# ---------------------------------------------------------------------------
function ConvertTo-ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = ConvertTo-ComRGB $officeThemeColors[$i - 1]
}
